# Horarios Linea 141 - actualizacion 04:52:25 (scrape refresh)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:52:25"
$ws1.Range("A3").Value = "Total filas: 19"

# rows 8-11: refresh scrape timestamp + recompute remaining minutes
$ws1.Range("A8").Value  = "04:52:25"
$ws1.Range("D8").Value  = 1
$ws1.Range("A9").Value  = "04:52:25"
$ws1.Range("D9").Value  = 24
$ws1.Range("A10").Value = "04:52:25"
$ws1.Range("D10").Value = 30
$ws1.Range("A11").Value = "04:52:25"
$ws1.Range("D11").Value = 42

# rows 14-22: refresh scrape timestamp + recompute remaining minutes
$ws1.Range("A14").Value = "04:52:25"
$ws1.Range("D14").Value = 54
$ws1.Range("A15").Value = "04:52:25"
$ws1.Range("D15").Value = 62
$ws1.Range("A16").Value = "04:52:25"
$ws1.Range("D16").Value = 72
$ws1.Range("A17").Value = "04:52:25"
$ws1.Range("D17").Value = 79
$ws1.Range("A18").Value = "04:52:25"
$ws1.Range("D18").Value = 82
$ws1.Range("A19").Value = "04:52:25"
$ws1.Range("D19").Value = 89
$ws1.Range("A20").Value = "04:52:25"
$ws1.Range("D20").Value = 95
$ws1.Range("A21").Value = "04:52:25"
$ws1.Range("D21").Value = 97
$ws1.Range("A22").Value = "04:52:25"
$ws1.Range("D22").Value = 99

# two new arrivals appended to the end of the schedule
$ws1.Range("A23").Value = "04:52:25"
$ws1.Range("B23").Value = "06:44"
$ws1.Range("C23").Value = "225_C ROCA-H SUR"
$ws1.Range("D23").Value = 112
$ws1.Range("E23").Value = "LP1912"

$ws1.Range("A24").Value = "04:52:25"
$ws1.Range("B24").Value = "06:46"
$ws1.Range("C24").Value = "215C_EL PATO"
$ws1.Range("D24").Value = 114
$ws1.Range("E24").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:52:25"
$ws2.Range("A3").Value = "Total filas: 5"

$ws2.Range("A7").Value = "04:52:25"
$ws2.Range("D7").Value = 42

$ws2.Range("A9").Value = "04:52:25"
$ws2.Range("D9").Value = 79

# new arrival appended to the end of the schedule
$ws2.Range("A10").Value = "04:52:25"
$ws2.Range("B10").Value = "06:46"
$ws2.Range("C10").Value = "215C_EL PATO"
$ws2.Range("D10").Value = 114
$ws2.Range("E10").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:52:25"

$ws3.Range("A6").Value = "04:52:25"
$ws3.Range("D6").Value = 52

$ws3.Range("A7").Value = "04:52:25"
$ws3.Range("D7").Value = 77

$ws3.Range("A8").Value = "04:52:25"
$ws3.Range("D8").Value = 101
